$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "63.673.78"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -2.10%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.129.52"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -2.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "607.72"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +1.42%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "145.65"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -5.13%  "
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.13%  "
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.125.49"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -2.22%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.520"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -1.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -3.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -3.84%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.470"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -1.81%  "
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000252"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -1.98%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.32"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -5.44%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.645.07"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.79%  "
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +2.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.705.99"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -2.43%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "3.133.82"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -2.65%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.82"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -3.25%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "475.33"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -2.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.55"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -2.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.713"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -1.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.81"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +0.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.54"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  -3.80%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "83.47"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.65%  "
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  +0.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -5.71%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.46"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -3.33%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.122"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -4.49%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "6.98"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +2.88%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -10.02%  "
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.68"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -3.03%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "26.17"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -3.30%  "
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +1.03%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0786"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  +6.60%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.93"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -3.73%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "52.87"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -4.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "457.21"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -3.86%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.01"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -8.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0392"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -3.69%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.118"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -6.39%  "
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "8.31"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -2.32%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.847.17"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -2.74%  "
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -6.47%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.266"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -5.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.42"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +1.59%  "
$cell.Style = "Normal"

$cell = $ws.Range("B48")
$cell.NumberFormat = "@"
$cell.Value = "USDe"
$cell.Style = "Normal"
$cell = $ws.Range("C48")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("B49")
$cell.NumberFormat = "@"
$cell.Value = "InjectiveProtocol"
$cell.Style = "Normal"
$cell = $ws.Range("C49")
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "26.22"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -4.72%  "
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -2.68%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "118.96"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -1.56%  "
$cell.Style = "Normal"
